# Got permission to publish from ADW!
#
# This script reproduces the edit described by the commit: several
# "RE-SCAN" / "AWAITING DUPLICATION" / "RRS" status notes are added to
# column D/H, and the strikethrough formatting that had been redacting a
# handful of dates in column E/A is removed now that publishing has been
# cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New status notes -------------------------------------------------
# Written in the same order the original authoring session created the
# shared-string table so new unique strings land in the expected slots:
# RE-SCAN, AWAITING DUPLICATION, "RE-SCAN, PAGE NUM", RRS, "RRS, PAGE NUM".
$ws.Range("H29").Value = "RE-SCAN"
$ws.Range("H31").Value = "RE-SCAN"
$ws.Range("H32").Value = "AWAITING DUPLICATION"
$ws.Range("H5").Value = "RE-SCAN, PAGE NUM"
$ws.Range("H4").Value = "RRS"
$ws.Range("D40").Value = "RRS"
$ws.Range("H48").Value = "RRS"
$ws.Range("D44").Value = "AWAITING DUPLICATION"
$ws.Range("D36").Value = "RRS, PAGE NUM"

# --- Remove strikethrough redaction now that publishing is cleared ----
$ws.Range("E5").Font.Strikethrough = $false
$ws.Range("E29").Font.Strikethrough = $false
$ws.Range("E31").Font.Strikethrough = $false
$ws.Range("E32").Font.Strikethrough = $false
$ws.Range("A44").Font.Strikethrough = $false

# --- Restore the cursor/selection left by the author -------------------
$ws.Range("B46").Select()
